# comment setLogger code, change files
#
# Remove the two runs ("1" and the long "ertfytry..." gibberish string)
# that precede the _GoBack bookmark in the first paragraph, leaving the
# bookmark start/end markers untouched.

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "1ertfytryipoiportiyopioprtiypoioprtiypoiportiypoiportiypoiproity",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 2)

Write-Output $found
